$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width (matches the new <cols> entry for column 1)
$ws.Columns.Item(1).ColumnWidth = 17.88671875

# Row 12: average of the k column (J), bold
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Rows 14-17: summary labels + stats
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the bold, size-12, vertically centred style once on B14, then
# propagate it to B15:B17 via a format copy/paste so every one of those
# cells lands on the same style index instead of the interpreter minting
# a fresh (unused) style per incremental property write.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Select A14:B17 and make A14 the active cell, matching the saved selection
$ws.Range("A14:B17").Select()

$wb.Save()
